# Fund ratios investor KYC test workbook - "fund 2" test data swap.
# Replaces the row 2/3 sample investor records with a second set of test
# fund identifiers, removes the extra sample rows (4-11) that belonged to
# the old "fund 1" bulk upload test, and tidies a couple of left-over
# blank styled cells near the end of row 2/3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: TSTF1 -> TSTF3 ------------------------------------------------
$ws.Range("I2").Value = "Kotak"
$ws.Range("C2").Value = "10/18/1991"

# --- Row 3: TSTF2 -> TSTF4 ------------------------------------------------
$ws.Range("I3").Value = "Axis"
$ws.Range("H3").Value = "Foreign"
$ws.Range("C3").Value = "05/30/2000"

# --- New fund identifiers (PAN + fund code) --------------------------------
$ws.Range("D2").Value = "TSTFU2121D"
$ws.Range("D3").Value = "TSTFU2222D"
$ws.Range("A2").Value = "TSTF3"
$ws.Range("B2").Value = "TSTF3"
$ws.Range("A3").Value = "TSTF4"
$ws.Range("B3").Value = "TSTF4"
# G2,G3 (Individual) / H2 (Domestic) / N2,O2,P2,N3,O3,P3 (Yes/No/No) stay as-is.

# --- Tidy the trailing blank-but-styled placeholder cells -----------------
# Drop the now-redundant Q2/R2 (and the matching Q3) blank cells entirely;
# S2 already carries the formatting that survives (it matched R2's style).
$ws.Range("Q2").Clear()
$ws.Range("R2").Clear()
$ws.Range("Q3").Clear()

# --- Drop the extra sample rows (old "fund 1" bulk rows 4-11) -------------
$ws.Range("A4:A11").EntireRow.Delete()

# --- Selection moves to B4 (now an empty row) ------------------------------
$ws.Range("B4").Select()
